$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.206045389175415
$ws.Range("B1").Value = 3.714653491973877
$ws.Range("C1").Value = 4.291790008544922
$ws.Range("D1").Value = 1.802934288978577
$ws.Range("E1").Value = 1.251164674758911
